# Org-chart personnel update:
#  - Slide 1 (Korean chart): "김동찬" -> "최성묵"
#  - Slide 2 (English chart):
#      * "Kim SonnChol, Oh Seung Hyun" -> two lines: "Kim SonnChol" / "Oh SeungHyun"
#      * "Joseph" -> "Cho Joseph"
#      * "Kim DongChan" -> "Choi " + "SeongMook" (new name, same role box)

$p = $ppt.ActivePresentation

# --- Slide 1: Korean org chart -------------------------------------------
$slide1 = $p.Slides.Item(1)

# Shape id=19 ("김동찬") -> "최성묵"
$slide1.Shapes.Item(12).TextFrame.TextRange.Text = "최성묵"

# --- Slide 2: English org chart ------------------------------------------
$slide2 = $p.Slides.Item(2)

# Shape id=15 ("Kim SonnChol, Oh Seung Hyun") -> split into two paragraphs
$tr15 = $slide2.Shapes.Item(8).TextFrame.TextRange
$tr15.Text = "Kim SonnChol"
$tr15.InsertAfter([char]13 + "Oh SeungHyun")

# Shape id=17 ("Joseph") -> "Cho Joseph"
$slide2.Shapes.Item(10).TextFrame.TextRange.Text = "Cho Joseph"

# Shape id=19 ("Kim DongChan") -> "Choi " + "SeongMook" (two runs, same paragraph)
$tr19 = $slide2.Shapes.Item(12).TextFrame.TextRange
$tr19.Text = "Choi "
$tr19.InsertAfter("SeongMook")
